$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Simple team-name swaps (shared string text for "NK Domzale" and
#     "NK Maribor" were swapped in the source data) ---
$ws.Range("F2").Value = "NK Maribor"
$ws.Range("F5").Value = "NK Domzale"
$ws.Range("F12").Value = "NK Maribor"
$ws.Range("E14").Value = "NK Domzale"
$ws.Range("F20").Value = "NK Domzale"
$ws.Range("E21").Value = "NK Maribor"
$ws.Range("E23").Value = "NK Domzale"
$ws.Range("F23").Value = "NK Maribor"
$ws.Range("E27").Value = "NK Maribor"
$ws.Range("F29").Value = "NK Domzale"
$ws.Range("F33").Value = "NK Domzale"
$ws.Range("F34").Value = "NK Maribor"
$ws.Range("E36").Value = "NK Maribor"
$ws.Range("E40").Value = "NK Domzale"
$ws.Range("E43").Value = "NK Maribor"
$ws.Range("E44").Value = "NK Domzale"
$ws.Range("F46").Value = "NK Maribor"
$ws.Range("F50").Value = "NK Domzale"
$ws.Range("E52").Value = "NK Domzale"
$ws.Range("E53").Value = "NK Maribor"
$ws.Range("E57").Value = "NK Maribor"
$ws.Range("F58").Value = "NK Domzale"
$ws.Range("F61").Value = "NK Maribor"
$ws.Range("E63").Value = "NK Domzale"
$ws.Range("E68").Value = "NK Maribor"
$ws.Range("F68").Value = "NK Domzale"
$ws.Range("F73").Value = "NK Maribor"
$ws.Range("E74").Value = "NK Domzale"
$ws.Range("F76").Value = "NK Maribor"
$ws.Range("F79").Value = "NK Domzale"
$ws.Range("E81").Value = "NK Maribor"
$ws.Range("E83").Value = "NK Domzale"
$ws.Range("F85").Value = "NK Maribor"
$ws.Range("F86").Value = "NK Domzale"
$ws.Range("F92").Value = "NK Maribor"
$ws.Range("F93").Value = "NK Domzale"
$ws.Range("E96").Value = "NK Maribor"
$ws.Range("E99").Value = "NK Domzale"
$ws.Range("F101").Value = "NK Maribor"
$ws.Range("F105").Value = "NK Maribor"
$ws.Range("E107").Value = "NK Domzale"
$ws.Range("F109").Value = "NK Domzale"
$ws.Range("E110").Value = "NK Maribor"
$ws.Range("E115").Value = "NK Domzale"
$ws.Range("F115").Value = "NK Maribor"
$ws.Range("F119").Value = "NK Domzale"
$ws.Range("F124").Value = "NK Maribor"
$ws.Range("E127").Value = "NK Maribor"
$ws.Range("E130").Value = "NK Domzale"
$ws.Range("E133").Value = "NK Domzale"
$ws.Range("E134").Value = "NK Maribor"
$ws.Range("F137").Value = "NK Domzale"
$ws.Range("F141").Value = "NK Maribor"
$ws.Range("E142").Value = "NK Domzale"
$ws.Range("E144").Value = "NK Maribor"
$ws.Range("F148").Value = "NK Domzale"
$ws.Range("E150").Value = "NK Maribor"
$ws.Range("F152").Value = "NK Domzale"
$ws.Range("E153").Value = "NK Maribor"
$ws.Range("E156").Value = "NK Domzale"
$ws.Range("F157").Value = "NK Maribor"
$ws.Range("E160").Value = "NK Maribor"
$ws.Range("F160").Value = "NK Domzale"
$ws.Range("F165").Value = "NK Maribor"
$ws.Range("E166").Value = "NK Domzale"
$ws.Range("E170").Value = "NK Maribor"
$ws.Range("E171").Value = "NK Domzale"

# --- Row 9 & Row 10: full match-data swap between the two fixtures ---
$ws.Range("B9").Value = 6814330
$ws.Range("F9").Value = "NK Aluminij"
$ws.Range("G9").Value = 1
$ws.Range("H9").Value = 0
$ws.Range("I9").Value = 0
$ws.Range("J9").Value = 0
$ws.Range("K9").Value = "H"
$ws.Range("L9").Value = 1.363
$ws.Range("M9").Value = 4.5
$ws.Range("N9").Value = 7
$ws.Range("O9").Value = 1.4
$ws.Range("P9").Value = 4.5
$ws.Range("Q9").Value = 7
$ws.Range("R9").Value = -1.25
$ws.Range("S9").Value = 1.85
$ws.Range("T9").Value = 1.95
$ws.Range("U9").Value = 2.75
$ws.Range("V9").Value = 1.8
$ws.Range("W9").Value = 2
$ws.Range("X9").Value = 0.3999999999999999
$ws.Range("Y9").Value = -1
$ws.Range("Z9").Value = -1
$ws.Range("AA9").Value = -0.5
$ws.Range("AB9").Value = 0.475
$ws.Range("AC9").Value = -1
$ws.Range("AD9").Value = 1
$ws.Range("B10").Value = 6814328
$ws.Range("F10").Value = "NK Bravo"
$ws.Range("G10").Value = 1
$ws.Range("H10").Value = 1
$ws.Range("I10").Value = 0
$ws.Range("J10").Value = 1
$ws.Range("K10").Value = "D"
$ws.Range("L10").Value = 2.35
$ws.Range("M10").Value = 3.1
$ws.Range("N10").Value = 2.9
$ws.Range("O10").Value = 2.15
$ws.Range("P10").Value = 3.1
$ws.Range("Q10").Value = 3.3
$ws.Range("R10").Value = -0.25
$ws.Range("S10").Value = 1.925
$ws.Range("T10").Value = 1.875
$ws.Range("U10").Value = 2.25
$ws.Range("V10").Value = 1.95
$ws.Range("W10").Value = 1.85
$ws.Range("X10").Value = -1
$ws.Range("Y10").Value = 2.1
$ws.Range("Z10").Value = -1
$ws.Range("AA10").Value = -0.5
$ws.Range("AB10").Value = 0.4375
$ws.Range("AC10").Value = -0.5
$ws.Range("AD10").Value = 0.425

# --- Row 174 & Row 176: full match-data swap between the two fixtures ---
$ws.Range("B174").Value = 7133777
$ws.Range("E174").Value = "NK Radomlje"
$ws.Range("F174").Value = "NK Celje"
$ws.Range("G174").Value = 1
$ws.Range("H174").Value = 1
$ws.Range("I174").Value = 1
$ws.Range("J174").Value = 0
$ws.Range("K174").Value = "D"
$ws.Range("L174").Value = 3.05
$ws.Range("M174").Value = 3.5
$ws.Range("N174").Value = 2
$ws.Range("O174").Value = 2.9
$ws.Range("P174").Value = 3.6
$ws.Range("Q174").Value = 2.1
$ws.Range("R174").Value = 0.25
$ws.Range("S174").Value = 1.9
$ws.Range("T174").Value = 1.9
$ws.Range("U174").Value = 2.75
$ws.Range("V174").Value = 1.8
$ws.Range("W174").Value = 2
$ws.Range("X174").Value = -1
$ws.Range("Y174").Value = 2.6
$ws.Range("Z174").Value = -1
$ws.Range("AA174").Value = 0.45
$ws.Range("AB174").Value = -0.5
$ws.Range("AC174").Value = -1
$ws.Range("AD174").Value = 1
$ws.Range("B176").Value = 7124153
$ws.Range("E176").Value = "NK Aluminij"
$ws.Range("F176").Value = "NK Domzale"
$ws.Range("G176").Value = 1
$ws.Range("H176").Value = 3
$ws.Range("I176").Value = 0
$ws.Range("J176").Value = 3
$ws.Range("K176").Value = "A"
$ws.Range("L176").Value = 2
$ws.Range("M176").Value = 3.6
$ws.Range("N176").Value = 3
$ws.Range("O176").Value = 1.333
$ws.Range("P176").Value = 4.75
$ws.Range("Q176").Value = 7
$ws.Range("R176").Value = -1.5
$ws.Range("S176").Value = 1.95
$ws.Range("T176").Value = 1.85
$ws.Range("U176").Value = 3.25
$ws.Range("V176").Value = 1.95
$ws.Range("W176").Value = 1.85
$ws.Range("X176").Value = -1
$ws.Range("Y176").Value = -1
$ws.Range("Z176").Value = 6
$ws.Range("AA176").Value = -1
$ws.Range("AB176").Value = 0.8500000000000001
$ws.Range("AC176").Value = 0.95
$ws.Range("AD176").Value = -1
